$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 15
$ws.Range("O2").Value = 1.2
$ws.Range("P2").Value = 4.5
$ws.Range("S2").Value = 2.63
$ws.Range("T2").Value = 1.5
# Row 3
$ws.Range("AP3").Value = 1.93
$ws.Range("AQ3").Value = 1.97
# Row 4
$ws.Range("S4").Value = 9
$ws.Range("T4").Value = 1.07
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.98
# Row 5
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 4
$ws.Range("L5").Value = 5
$ws.Range("Z5").Value = 8
$ws.Range("AB5").Value = 17
$ws.Range("AJ5").Value = 8.5
# Row 9
$ws.Range("G9").Value = 2.8
$ws.Range("I9").Value = 2.35
$ws.Range("J9").Value = 3.6
$ws.Range("L9").Value = 3.2
$ws.Range("AJ9").Value = 7.5
$ws.Range("AL9").Value = 10
$ws.Range("AM9").Value = 23
# Row 15
$ws.Range("G15").Value = 1.36
$ws.Range("I15").Value = 7.5
$ws.Range("J15").Value = 1.83
$ws.Range("O15").Value = 1.17
$ws.Range("P15").Value = 5
$ws.Range("Q15").Value = 1.57
$ws.Range("R15").Value = 2.35
$ws.Range("S15").Value = 2.38
$ws.Range("T15").Value = 1.53
$ws.Range("U15").Value = 1.29
$ws.Range("V15").Value = 3.5
$ws.Range("AB15").Value = 9
$ws.Range("AJ15").Value = 21
$ws.Range("AR15").Value = 1.98
$ws.Range("AS15").Value = 1.88
# Row 17
$ws.Range("G17").Value = 3.9
$ws.Range("K17").Value = 2.63
$ws.Range("M17").Value = 1.01
$ws.Range("N17").Value = 26
$ws.Range("O17").Value = 1.1
$ws.Range("P17").Value = 7
$ws.Range("Q17").Value = 1.36
$ws.Range("R17").Value = 3.1
$ws.Range("S17").Value = 1.83
$ws.Range("T17").Value = 1.83
$ws.Range("W17").Value = 1.4
$ws.Range("X17").Value = 2.75
$ws.Range("Y17").Value = 21
$ws.Range("AA17").Value = 13
$ws.Range("AD17").Value = 23
$ws.Range("AE17").Value = 26
$ws.Range("AF17").Value = 9.5
$ws.Range("AK17").Value = 13
# Row 18
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 5.5
$ws.Range("I18").Value = 1.33
$ws.Range("J18").Value = 6.5
$ws.Range("K18").Value = 2.88
$ws.Range("L18").Value = 1.73
$ws.Range("N18").Value = 26
$ws.Range("Q18").Value = 1.33
$ws.Range("R18").Value = 3.4
$ws.Range("U18").Value = 1.2
$ws.Range("V18").Value = 4.33
$ws.Range("W18").Value = 1.57
$ws.Range("X18").Value = 2.25
$ws.Range("Y18").Value = 29
$ws.Range("AA18").Value = 21
$ws.Range("AB18").Value = 81
$ws.Range("AD18").Value = 41
$ws.Range("AE18").Value = 26
$ws.Range("AF18").Value = 12
$ws.Range("AG18").Value = 17
$ws.Range("AI18").Value = 126
$ws.Range("AK18").Value = 9.5
$ws.Range("AM18").Value = 11
$ws.Range("AN18").Value = 10
$ws.Range("AO18").Value = 19
# Row 19
$ws.Range("G19").Value = 2.9
$ws.Range("I19").Value = 2.25
$ws.Range("J19").Value = 3.25
$ws.Range("L19").Value = 2.75
$ws.Range("Y19").Value = 15
$ws.Range("Z19").Value = 19
$ws.Range("AC19").Value = 21
$ws.Range("AD19").Value = 23
$ws.Range("AJ19").Value = 12
$ws.Range("AL19").Value = 9.5
$ws.Range("AM19").Value = 23
$ws.Range("AN19").Value = 15
# Row 20
$ws.Range("M20").Value = 1.05
$ws.Range("N20").Value = 11
$ws.Range("AP20").Value = 1.43
$ws.Range("AQ20").Value = 2.85
# Row 21
$ws.Range("G21").Value = 3.1
$ws.Range("H21").Value = 2.8
$ws.Range("I21").Value = 2.63
$ws.Range("J21").Value = 3.75
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 3.25
$ws.Range("M21").Value = 1.1
$ws.Range("N21").Value = 7
$ws.Range("O21").Value = 1.4
$ws.Range("P21").Value = 2.75
$ws.Range("Q21").Value = 2.35
$ws.Range("R21").Value = 1.57
$ws.Range("S21").Value = 4.33
$ws.Range("T21").Value = 1.2
$ws.Range("U21").Value = 1.5
$ws.Range("V21").Value = 2.5
$ws.Range("W21").Value = 1.95
$ws.Range("X21").Value = 1.8
$ws.Range("Y21").Value = 8.5
$ws.Range("AA21").Value = 12
$ws.Range("AC21").Value = 26
$ws.Range("AD21").Value = 41
$ws.Range("AE21").Value = 7
$ws.Range("AG21").Value = 15
$ws.Range("AH21").Value = 51
$ws.Range("AI21").Value = 301
$ws.Range("AJ21").Value = 7.5
$ws.Range("AK21").Value = 12
$ws.Range("AL21").Value = 11
$ws.Range("AM21").Value = 26
$ws.Range("AN21").Value = 23
$ws.Range("AO21").Value = 34
# Row 22
$ws.Range("I22").Value = 3.75
$ws.Range("M22").Value = 1.06
$ws.Range("N22").Value = 10
$ws.Range("Q22").Value = 2.08
$ws.Range("R22").Value = 1.73
# Row 23
$ws.Range("G23").Value = 3.6
$ws.Range("H23").Value = 3.2
$ws.Range("I23").Value = 1.91
$ws.Range("J23").Value = 4.5
$ws.Range("K23").Value = 2.05
$ws.Range("L23").Value = 2.75
$ws.Range("Y23").Value = 9.5
$ws.Range("Z23").Value = 19
$ws.Range("AA23").Value = 13
$ws.Range("AC23").Value = 34
$ws.Range("AK23").Value = 8.5
$ws.Range("AL23").Value = 9
$ws.Range("AM23").Value = 17
$ws.Range("AN23").Value = 17
# Row 24
$ws.Range("I24").Value = 5
$ws.Range("K24").Value = 2.05
$ws.Range("N24").Value = 7.5
$ws.Range("U24").Value = 1.5
$ws.Range("V24").Value = 2.5
$ws.Range("Y24").Value = 5.5
$ws.Range("Z24").Value = 7
$ws.Range("AN24").Value = 41
# Row 26
$ws.Range("G26").Value = 3.25
$ws.Range("H26").Value = 2.9
$ws.Range("I26").Value = 2.35
$ws.Range("L26").Value = 3.2
$ws.Range("O26").Value = 1.5
$ws.Range("P26").Value = 2.5
$ws.Range("Q26").Value = 2.5
$ws.Range("R26").Value = 1.5
$ws.Range("S26").Value = 5
$ws.Range("T26").Value = 1.17
$ws.Range("U26").Value = 1.57
$ws.Range("V26").Value = 2.25
$ws.Range("W26").Value = 2.1
$ws.Range("X26").Value = 1.67
$ws.Range("Y26").Value = 8
$ws.Range("AE26").Value = 6.5
$ws.Range("AH26").Value = 67
$ws.Range("AM26").Value = 23
$ws.Range("AN26").Value = 23
$ws.Range("AO26").Value = 41
$ws.Range("AP26").Value = 1.93
$ws.Range("AQ26").Value = 1.93
# Row 27
$ws.Range("G27").Value = 1.85
$ws.Range("H27").Value = 3.6
$ws.Range("I27").Value = 3.9
$ws.Range("J27").Value = 2.5
$ws.Range("M27").Value = 1.03
$ws.Range("N27").Value = 15
$ws.Range("O27").Value = 1.18
$ws.Range("P27").Value = 4.5
$ws.Range("Q27").Value = 1.65
$ws.Range("R27").Value = 2.2
$ws.Range("S27").Value = 2.5
$ws.Range("T27").Value = 1.5
$ws.Range("W27").Value = 1.62
$ws.Range("X27").Value = 2.2
$ws.Range("Z27").Value = 10
$ws.Range("AA27").Value = 8.5
$ws.Range("AC27").Value = 13
$ws.Range("AG27").Value = 13
$ws.Range("AJ27").Value = 15
$ws.Range("AN27").Value = 29
